# Updated cryptos list on Sun Mar 17 07:15:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = "65.146.38";   E = "  -5.92%  " }
    3  = @{ D = "3.450.61";    E = "  -7.45%  " }
    4  = @{               E = "  +0.14%  " }
    5  = @{ D = "560.02";      E = "  -8.52%  " }
    6  = @{ D = "182.38";      E = "  -4.88%  " }
    7  = @{ D = "3.447.10";    E = "  -7.41%  " }
    8  = @{ D = "0.596";       E = "  -6.58%  " }
    9  = @{               E = "  +0.17%  " }
    10 = @{ D = "0.642";       E = "  -11.49%  " }
    11 = @{               E = "  -12.95%  " }
    12 = @{ D = "51.05";       E = "  -15.27%  " }
    13 = @{ D = "0.0000249";   E = "  -14.58%  " }
    14 = @{ D = "9.41";        E = "  -11.68%  " }
    15 = @{ D = "4.005.33";    E = "  -7.18%  " }
    16 = @{               E = "  -1.63%  " }
    17 = @{ D = "3.457.24";    E = "  -7.06%  " }
    18 = @{ B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "64.911.76"; E = "  -6.02%  " }
    19 = @{ B = "Chainlink";  C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "17.71";     E = "  -8.84%  " }
    20 = @{ D = "11.67";       E = "  -10.00%  " }
    21 = @{ D = "1.02";        E = "  -10.66%  " }
    22 = @{ D = "378.32";      E = "  -8.42%  " }
    23 = @{ D = "4.12";        E = "  -10.40%  " }
    24 = @{ D = "82.89";       E = "  -7.61%  " }
    25 = @{ D = "10.62";       E = "  -2.67%  " }
    26 = @{ D = "2.79";        E = "  -8.92%  " }
    27 = @{               E = "  -0.91%  " }
    28 = @{ D = "11.87";       E = "  -7.91%  " }
    29 = @{ D = "3.38";        E = "  -11.45%  " }
    30 = @{ D = "8.50";        E = "  -12.62%  " }
    31 = @{ D = "30.04";       E = "  -9.54%  " }
    32 = @{ D = "6.94";        E = "  -9.65%  " }
    33 = @{ D = "602.28";      E = "  -4.85%  " }
    34 = @{ D = "11.76";       E = "  -8.31%  " }
    35 = @{ D = "62.36";       E = "  -5.09%  " }
    36 = @{ D = "0.109";       E = "  -11.78%  " }
    37 = @{ D = "40.25";       E = "  -12.58%  " }
    38 = @{ D = "0.999";       E = "  -0.08%  " }
    39 = @{ D = "0.387";       E = "  -6.60%  " }
    40 = @{ D = "0.0₃0716";    E = "  -14.14%  " }
    41 = @{               E = "  +0.06%  " }
    42 = @{               E = "  -9.56%  " }
    43 = @{ D = "2.915.47";    E = "  +1.12%  " }
    44 = @{               E = "  -11.60%  " }
    45 = @{               E = "  -8.73%  " }
    46 = @{ B = "VeChain";      C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";   D = "0.0391"; E = "  -12.70%  " }
    47 = @{ B = "ApeXProtocol"; C = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex";  D = "3.05";    E = "  -2.81%  " }
    48 = @{ D = "0.126";       E = "  -9.92%  " }
    49 = @{ D = "136.75";      E = "  -4.01%  " }
    50 = @{ D = "2.45";        E = "  -11.06%  " }
    51 = @{ D = "8.14";        E = "  -11.66%  " }
}

# Price values in column D are textual (e.g. "65.146.38", "3.05") and must
# stay stored as text, exactly as they were before the edit (inline strings).
# Values that look like plain numbers would otherwise get silently coerced
# to floating point numbers by Excel, so force those particular cells to
# text format before writing them.
foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    if ($rowData.ContainsKey("D")) {
        $dValue = $rowData["D"]
        $isNumeric = $dValue -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'
        if ($isNumeric) {
            $ws.Range("D$row").NumberFormat = "@"
        }
    }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $rowData[$col]
    }
}
